$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 2).Value = 1379508
$ws.Cells.Item(4, 3).Value = 11870
$ws.Cells.Item(4, 4).Value = 259073
$ws.Cells.Item(4, 5).Value = 1039150
$ws.Cells.Item(4, 7).Value = 498
$ws.Cells.Item(4, 8).Value = 81285

$ws.Cells.Item(10, 6).Value = 1576

$ws.Cells.Item(17, 2).Value = 68822
$ws.Cells.Item(17, 3).Value = 1515
$ws.Cells.Item(17, 4).Value = 22406
$ws.Cells.Item(17, 5).Value = 44455
$ws.Cells.Item(17, 6).Value = 785
$ws.Cells.Item(17, 7).Value = 72
$ws.Cells.Item(17, 8).Value = 1961

$ws.Cells.Item(66, 4).Value = 1250
$ws.Cells.Item(66, 5).Value = 2306
$ws.Cells.Item(66, 6).Value = 27

$ws.Cells.Item(108, 2).Value = 801
$ws.Cells.Item(108, 3).Value = 9
$ws.Cells.Item(108, 4).Value = 517
$ws.Cells.Item(108, 5).Value = 277

$ws.Cells.Item(173, 4).Value = 24
$ws.Cells.Item(173, 5).Value = 30

$ws.Cells.Item(174, 2).Value = 56
$ws.Cells.Item(174, 3).Value = 5
$ws.Cells.Item(174, 5).Value = 46
$ws.Cells.Item(174, 7).Value = 1
$ws.Cells.Item(174, 8).Value = 9

# Row 187 currently Granada, becomes Gambia with updated stats
$ws.Cells.Item(187, 1).Value = "Gambia"
$ws.Cells.Item(187, 2).Value = 22
$ws.Cells.Item(187, 3).Value = 2
$ws.Cells.Item(187, 4).Value = 10
$ws.Cells.Item(187, 5).Value = 11
$ws.Cells.Item(187, 6).Value = 0
$ws.Cells.Item(187, 8).Value = 1

# Row 188 currently Gambia, becomes Granada with updated stats
$ws.Cells.Item(188, 1).Value = "Granada"
$ws.Cells.Item(188, 2).Value = 21
$ws.Cells.Item(188, 4).Value = 13
$ws.Cells.Item(188, 5).Value = 8
$ws.Cells.Item(188, 6).Value = 4
$ws.Cells.Item(188, 8).Value = 0

# Row 192 currently Belice, becomes Nueva Caledonia with updated stats
$ws.Cells.Item(192, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(192, 4).Value = 18
$ws.Cells.Item(192, 8).Value = 0

# Row 193 currently Nueva Caledonia, becomes Belice with updated stats
$ws.Cells.Item(193, 1).Value = "Belice"
$ws.Cells.Item(193, 4).Value = 16
$ws.Cells.Item(193, 8).Value = 2
